$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.665.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.517.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.38%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.22"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.903.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.493.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.862"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.691.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.09%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.87"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.09"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.66"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.36"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.34"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.97%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.12%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.18"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.18%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0301"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.036.67"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.39"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.81"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.761.24"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.05%  "
